# Auto-generated Excel COM-interop script applying the Seraph_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1074.9348
$ws.Range("I15").Value = 1074.9348
$ws.Range("K15").Value = 3224.8044
$ws.Range("M15").Value = -3055.8044
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
# Row 135
$ws.Range("H135").Value = 1368.2
$ws.Range("I135").Value = 1212.7727
$ws.Range("K135").Value = 10914.9543
$ws.Range("M135").Value = -8379.9543
# Row 137
$ws.Range("H137").Value = 2643
$ws.Range("I137").Value = 2043.7059
$ws.Range("K137").Value = 6131.1177
$ws.Range("M137").Value = -3581.1177
# Row 138
$ws.Range("H138").Value = 3656.2144
$ws.Range("I138").Value = 9363
$ws.Range("J138").Value = 2705.0833
$ws.Range("K138").Value = 28089
$ws.Range("L138").Value = 8115.249899999999
$ws.Range("M138").Value = -22949
$ws.Range("N138").Value = -18395.2499

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6507.7935
$ws.Range("I32").Value = 4117.4526
$ws.Range("J32").Value = 19176.6
$ws.Range("K32").Value = 4117.4526
$ws.Range("L32").Value = 19176.6
$ws.Range("M32").Value = -3830.4526
$ws.Range("N32").Value = -19750.6
# Row 61
$ws.Range("H61").Value = 1628.1305
$ws.Range("I61").Value = 1579.4546
$ws.Range("J61").Value = 2699
$ws.Range("K61").Value = 1579.4546
$ws.Range("L61").Value = 2699
$ws.Range("M61").Value = -1367.4546
$ws.Range("N61").Value = -3123
# Row 74
$ws.Range("H74").Value = 2006.6538
$ws.Range("I74").Value = 1133.9546
$ws.Range("J74").Value = 6806.5
$ws.Range("K74").Value = 1133.9546
$ws.Range("L74").Value = 6806.5
$ws.Range("M74").Value = -259.9546
$ws.Range("N74").Value = -8554.5
# Row 77
$ws.Range("H77").Value = 2006.6538
$ws.Range("I77").Value = 1133.9546
$ws.Range("J77").Value = 6806.5
$ws.Range("K77").Value = 5669.773
$ws.Range("L77").Value = 34032.5
$ws.Range("M77").Value = -1301.773
$ws.Range("N77").Value = -42768.5
# Row 136
$ws.Range("H136").Value = 1628.1305
$ws.Range("I136").Value = 1579.4546
$ws.Range("J136").Value = 2699
$ws.Range("K136").Value = 4738.3638
$ws.Range("L136").Value = 8097
$ws.Range("M136").Value = -2188.3638
$ws.Range("N136").Value = -13197

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 3385.9524
$ws.Range("I58").Value = 1220.1111
$ws.Range("K58").Value = 1220.1111
$ws.Range("M58").Value = -1017.1111
# Row 131
$ws.Range("H131").Value = 80000
$ws.Range("J131").Value = 80000
$ws.Range("L131").Value = 80000
$ws.Range("N131").Value = -90080
# Row 132
$ws.Range("H132").Value = 2207.5356
$ws.Range("I132").Value = 2287.92
$ws.Range("J132").Value = 1537.6666
$ws.Range("K132").Value = 6863.76
$ws.Range("L132").Value = 4612.9998
$ws.Range("M132").Value = -4333.76
$ws.Range("N132").Value = -9672.9998
# Row 136
$ws.Range("H136").Value = 3385.9524
$ws.Range("I136").Value = 1220.1111
$ws.Range("K136").Value = 3660.3333
$ws.Range("M136").Value = -1110.3333

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 266.33334
$ws.Range("I23").Value = 199.5
$ws.Range("K23").Value = 598.5
$ws.Range("M23").Value = -363.5
# Row 43
$ws.Range("H43").Value = 450
$ws.Range("I43").Value = 450
$ws.Range("K43").Value = 1350
$ws.Range("M43").Value = -1236
# Row 137
$ws.Range("H137").Value = 7126
$ws.Range("I137").Value = 4712.5
$ws.Range("K137").Value = 14137.5
$ws.Range("M137").Value = -9037.5

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 1765001
$ws.Range("I11").Value = 5000000
$ws.Range("J11").Value = 686668
$ws.Range("K11").Value = 5000000
$ws.Range("L11").Value = 686668
$ws.Range("M11").Value = -4999861
$ws.Range("N11").Value = -686946
# Row 102
$ws.Range("H102").Value = 2272.7058
$ws.Range("I102").Value = 927.8333
$ws.Range("K102").Value = 927.8333
$ws.Range("M102").Value = 694.1667
# Row 132
$ws.Range("H132").Value = 4148.125
$ws.Range("I132").Value = 1933
$ws.Range("J132").Value = 5477.2
$ws.Range("K132").Value = 5799
$ws.Range("L132").Value = 16431.6
$ws.Range("M132").Value = -3269
$ws.Range("N132").Value = -21491.6

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 600
$ws.Range("I7").Value = 400
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 400
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = -288
$ws.Range("N7").Value = -1224
# Row 40
$ws.Range("H40").Value = 4529.143
$ws.Range("I40").Value = 4176
$ws.Range("K40").Value = 4176
$ws.Range("M40").Value = -4040
# Row 68
$ws.Range("H68").Value = 3439.6667
$ws.Range("I68").Value = 3122
$ws.Range("J68").Value = 4075
$ws.Range("K68").Value = 3122
$ws.Range("L68").Value = 4075
$ws.Range("M68").Value = -2373
$ws.Range("N68").Value = -5573
# Row 71
$ws.Range("H71").Value = 3439.6667
$ws.Range("I71").Value = 3122
$ws.Range("J71").Value = 4075
$ws.Range("K71").Value = 15610
$ws.Range("L71").Value = 20375
$ws.Range("M71").Value = -11866
$ws.Range("N71").Value = -27863
# Row 122
$ws.Range("H122").Value = 4553.8945
$ws.Range("I122").Value = 3215.4285
$ws.Range("K122").Value = 9646.2855
$ws.Range("M122").Value = -7196.2855
# Row 126
$ws.Range("H126").Value = 600
$ws.Range("I126").Value = 400
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 1200
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = 1270
$ws.Range("N126").Value = -7940
# Row 132
$ws.Range("H132").Value = 4844.96
$ws.Range("I132").Value = 4289.231
$ws.Range("K132").Value = 12867.693
$ws.Range("M132").Value = -10337.693

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
# Row 81
$ws.Range("H81").Value = 2667.6
$ws.Range("I81").Value = 3178.375
$ws.Range("J81").Value = 624.5
$ws.Range("K81").Value = 6356.75
$ws.Range("L81").Value = 1249
$ws.Range("M81").Value = -5295.75
$ws.Range("N81").Value = -3371
# Row 84
$ws.Range("H84").Value = 2667.6
$ws.Range("I84").Value = 3178.375
$ws.Range("J84").Value = 624.5
$ws.Range("K84").Value = 31783.75
$ws.Range("L84").Value = 6245
$ws.Range("M84").Value = -26479.75
$ws.Range("N84").Value = -16853
# Row 122
$ws.Range("H122").Value = 1644.1333
$ws.Range("I122").Value = 1690.3572
$ws.Range("J122").Value = 997
$ws.Range("K122").Value = 5071.071599999999
$ws.Range("L122").Value = 2991
$ws.Range("M122").Value = -2621.071599999999
$ws.Range("N122").Value = -7891
# Row 132
$ws.Range("H132").Value = 603.2308
$ws.Range("I132").Value = 603.2308
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1809.6924
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 720.3075999999999
$ws.Range("N132").ClearContents()
